$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite header row (A1:L1) with the new column set / order
$ws.Range("A1").Value = "first name"
$ws.Range("B1").Value = "last name"
$ws.Range("C1").Value = "channel"
$ws.Range("D1").Value = "email"
$ws.Range("E1").Value = "phone number"
$ws.Range("F1").Value = "Source"
$ws.Range("G1").Value = "Position"
$ws.Range("H1").Value = "company"
$ws.Range("I1").Value = "Notes"
$ws.Range("J1").Value = "referral email"
$ws.Range("K1").Value = "gender"
$ws.Range("L1").Value = "date"

# New custom column widths for columns I (9) and J (10)
$ws.Columns.Item(9).ColumnWidth = 31.5
$ws.Columns.Item(10).ColumnWidth = 12.5

# Page orientation set to portrait (page setup touched during local conversion)
$ws.PageSetup.Orientation = 1

# Selection moved to I15
$ws.Range("I15").Select()
